# New weekly price record for "Perejil" @ Vega Modelo de Temuco.
# The source table is sorted by Fecha (column D) descending, so the newest
# observation (2023-11-28) is inserted at the top of the data block (row 465,
# right after the header row), pushing every existing record down by one row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 465:end down by one, carrying formatting (date style, etc.) down
# with them, and leave row 465 blank/ready for the new record.
$ws.Rows("465:465").Insert()

# Populate the newly inserted row with the new observation.
$ws.Range("A465").Value = 10
$ws.Range("B465").Value = "Vega Modelo de Temuco"
$ws.Range("C465").Value = "La Araucanía"
$ws.Range("D465").Value = "2023-11-28"
$ws.Range("E465").Value = 9
$ws.Range("F465").Value = 100112044
$ws.Range("G465").Value = "Perejil"
$ws.Range("H465").Value = "Sin especificar"
$ws.Range("I465").Value = "Primera"
$ws.Range("J465").Value = 45
$ws.Range("K465").Value = 5000
$ws.Range("L465").Value = 5000
$ws.Range("M465").Value = 5000
$ws.Range("N465").Value = "$/docena de atados (3 kilos)"
$ws.Range("O465").Value = "Provincia de Cautín"
$ws.Range("P465").Value = 1667
$ws.Range("Q465").Value = 3
$ws.Range("R465").Value = "Hortaliza"
